# Auto-generated Excel COM-interop script to apply market data refresh
# (currentAveragePrice / LevePrice / LeveProfit columns) across sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 2250.238
$ws.Range("I125").Value = 2639
$ws.Range("J125").Value = 2128.75
$ws.Range("K125").Value = 23751
$ws.Range("L125").Value = 19158.75
$ws.Range("M125").Value = -21291
$ws.Range("N125").Value = -24078.75

$ws.Range("H132").Value = 3626800
$ws.Range("I132").Value = 3971976.2
$ws.Range("J132").Value = 2449.3333
$ws.Range("K132").Value = 11915928.6
$ws.Range("L132").Value = 7347.999899999999
$ws.Range("M132").Value = -11913398.6
$ws.Range("N132").Value = -12407.9999

$ws.Range("H137").Value = 1469.3846
$ws.Range("I137").Value = 1646.4117
$ws.Range("J137").Value = 1135
$ws.Range("K137").Value = 4939.2351
$ws.Range("L137").Value = 3405
$ws.Range("M137").Value = -2389.2351
$ws.Range("N137").Value = -8505

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7247.3057
$ws.Range("I32").Value = 6188.1177
$ws.Range("K32").Value = 6188.1177
$ws.Range("M32").Value = -5901.1177

$ws.Range("H45").Value = 78613.766
$ws.Range("I45").Value = 112621.11
$ws.Range("J45").Value = 2097.25
$ws.Range("K45").Value = 112621.11
$ws.Range("L45").Value = 2097.25
$ws.Range("M45").Value = -112244.11
$ws.Range("N45").Value = -2851.25

$ws.Range("H97").Value = 29621.371
$ws.Range("I97").Value = 39148.42
$ws.Range("J97").Value = 2098.7778
$ws.Range("K97").Value = 39148.42
$ws.Range("L97").Value = 2098.7778
$ws.Range("M97").Value = -38652.42
$ws.Range("N97").Value = -3090.7778

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 47739.047
$ws.Range("I20").Value = 69344.266
$ws.Range("J20").Value = 1442.1428
$ws.Range("K20").Value = 69344.266
$ws.Range("L20").Value = 1442.1428
$ws.Range("M20").Value = -69097.266
$ws.Range("N20").Value = -1936.1428

$ws.Range("H31").Value = 23026
$ws.Range("J31").Value = 23026
$ws.Range("L31").Value = 23026
$ws.Range("N31").Value = -23530

$ws.Range("H99").Value = 1609.8182
$ws.Range("I99").Value = 1106
$ws.Range("J99").Value = 1981.0526
$ws.Range("K99").Value = 1106
$ws.Range("L99").Value = 1981.0526
$ws.Range("M99").Value = 392
$ws.Range("N99").Value = -4977.0526

$ws.Range("H105").Value = 49249.855
$ws.Range("I105").Value = 73123.36
$ws.Range("J105").Value = 1502.8572
$ws.Range("K105").Value = 73123.36
$ws.Range("L105").Value = 1502.8572
$ws.Range("M105").Value = -71376.36
$ws.Range("N105").Value = -4996.8572

$ws.Range("H107").Value = 71462350
$ws.Range("I107").Value = 100046420
$ws.Range("J107").Value = 2203
$ws.Range("K107").Value = 100046420
$ws.Range("L107").Value = 2203
$ws.Range("M107").Value = -100044500
$ws.Range("N107").Value = -6043

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 10666.667
$ws.Range("I17").Value = 3000
$ws.Range("J17").Value = 14500
$ws.Range("K17").Value = 3000
$ws.Range("L17").Value = 14500
$ws.Range("M17").Value = -2826
$ws.Range("N17").Value = -14848

$ws.Range("H58").Value = 2680.8333
$ws.Range("I58").Value = 1871.3636
$ws.Range("J58").Value = 3952.8572
$ws.Range("K58").Value = 1871.3636
$ws.Range("L58").Value = 3952.8572
$ws.Range("M58").Value = -1668.3636
$ws.Range("N58").Value = -4358.8572

$ws.Range("H107").Value = 6993.9414
$ws.Range("I107").Value = 8815.154
$ws.Range("J107").Value = 1075
$ws.Range("K107").Value = 8815.154
$ws.Range("L107").Value = 1075
$ws.Range("M107").Value = -6895.154
$ws.Range("N107").Value = -4915

$ws.Range("H136").Value = 2680.8333
$ws.Range("I136").Value = 1871.3636
$ws.Range("J136").Value = 3952.8572
$ws.Range("K136").Value = 5614.0908
$ws.Range("L136").Value = 11858.5716
$ws.Range("M136").Value = -3064.0908
$ws.Range("N136").Value = -16958.5716

$ws.Range("H141").Value = 83434.45
$ws.Range("J141").Value = 58472.375
$ws.Range("L141").Value = 58472.375
$ws.Range("N141").Value = -68832.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H59").Value = 875.8333
$ws.Range("I59").Value = 727.5
$ws.Range("J59").Value = 950
$ws.Range("K59").Value = 2182.5
$ws.Range("L59").Value = 2850
$ws.Range("M59").Value = -1642.5
$ws.Range("N59").Value = -3930

$ws.Range("H114").Value = 898.5
$ws.Range("I114").Value = 472
$ws.Range("J114").Value = 3031
$ws.Range("K114").Value = 1416
$ws.Range("L114").Value = 9093
$ws.Range("M114").Value = 1838
$ws.Range("N114").Value = -15601

$ws.Range("H132").Value = 1687.4706
$ws.Range("I132").Value = 846.4
$ws.Range("J132").Value = 2037.9166
$ws.Range("K132").Value = 7617.599999999999
$ws.Range("L132").Value = 18341.2494
$ws.Range("M132").Value = -5087.599999999999
$ws.Range("N132").Value = -23401.2494

$ws.Range("H137").Value = 52013.773
$ws.Range("I137").Value = 101401
$ws.Range("J137").Value = 10857.75
$ws.Range("K137").Value = 304203
$ws.Range("L137").Value = 32573.25
$ws.Range("M137").Value = -299103
$ws.Range("N137").Value = -42773.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 70006
$ws.Range("J19").Value = 70006
$ws.Range("L19").Value = 70006
$ws.Range("N19").Value = -70582

$ws.Range("H107").Value = 631705.2
$ws.Range("I107").Value = 403.14285
$ws.Range("J107").Value = 1122717.9
$ws.Range("K107").Value = 403.14285
$ws.Range("L107").Value = 1122717.9
$ws.Range("M107").Value = 1516.85715
$ws.Range("N107").Value = -1126557.9

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 9799
$ws.Range("J5").Value = 9799
$ws.Range("L5").Value = 9799
$ws.Range("N5").Value = -10025

$ws.Range("H7").Value = 3656.0588
$ws.Range("I7").Value = 3930.9443
$ws.Range("J7").Value = 3346.8125
$ws.Range("K7").Value = 3930.9443
$ws.Range("L7").Value = 3346.8125
$ws.Range("M7").Value = -3818.9443
$ws.Range("N7").Value = -3570.8125

$ws.Range("H40").Value = 57478.832
$ws.Range("I40").Value = 101461.9
$ws.Range("K40").Value = 101461.9
$ws.Range("M40").Value = -101325.9

$ws.Range("H46").Value = 844020.75
$ws.Range("I46").Value = 496.66666
$ws.Range("J46").Value = 1125195.5
$ws.Range("K46").Value = 496.66666
$ws.Range("L46").Value = 1125195.5
$ws.Range("M46").Value = -308.66666
$ws.Range("N46").Value = -1125571.5

$ws.Range("H61").Value = 1919
$ws.Range("I61").Value = 2051.75
$ws.Range("J61").Value = 1786.25
$ws.Range("K61").Value = 2051.75
$ws.Range("L61").Value = 1786.25
$ws.Range("M61").Value = -1849.75
$ws.Range("N61").Value = -2190.25

$ws.Range("H100").Value = 2281.2856
$ws.Range("I100").Value = 2066.6667
$ws.Range("J100").Value = 2442.25
$ws.Range("K100").Value = 2066.6667
$ws.Range("L100").Value = 2442.25
$ws.Range("M100").Value = -1525.6667
$ws.Range("N100").Value = -3524.25

$ws.Range("H113").Value = 1919
$ws.Range("I113").Value = 2051.75
$ws.Range("J113").Value = 1786.25
$ws.Range("K113").Value = 2051.75
$ws.Range("L113").Value = 1786.25
$ws.Range("M113").Value = 118.25
$ws.Range("N113").Value = -6126.25

$ws.Range("H122").Value = 2345
$ws.Range("I122").Value = 2345
$ws.Range("K122").Value = 7035
$ws.Range("M122").Value = -4585

$ws.Range("H126").Value = 3656.0588
$ws.Range("I126").Value = 3930.9443
$ws.Range("J126").Value = 3346.8125
$ws.Range("K126").Value = 11792.8329
$ws.Range("L126").Value = 10040.4375
$ws.Range("M126").Value = -9322.832900000001
$ws.Range("N126").Value = -14980.4375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 8772.091
$ws.Range("J15").Value = 8909.299999999999
$ws.Range("L15").Value = 8909.299999999999
$ws.Range("N15").Value = -9485.299999999999

$ws.Range("H107").Value = 154854.08
$ws.Range("I107").Value = 1067.1111
$ws.Range("J107").Value = 500874.75
$ws.Range("K107").Value = 3201.3333
$ws.Range("L107").Value = 1502624.25
$ws.Range("M107").Value = -1281.3333
$ws.Range("N107").Value = -1506464.25

$ws.Range("H122").Value = 2470.6365
$ws.Range("I122").Value = 1951
$ws.Range("J122").Value = 2767.5715
$ws.Range("K122").Value = 5853
$ws.Range("L122").Value = 8302.7145
$ws.Range("M122").Value = -3403
$ws.Range("N122").Value = -13202.7145

$ws.Range("H126").Value = 2540.2856
$ws.Range("I126").Value = 2456.8
$ws.Range("K126").Value = 7370.400000000001
$ws.Range("M126").Value = -4900.400000000001

$ws.Range("H136").Value = 1443.1864
$ws.Range("I136").Value = 522.65515
$ws.Range("J136").Value = 2333.0334
$ws.Range("K136").Value = 1567.96545
$ws.Range("L136").Value = 6999.100199999999
$ws.Range("M136").Value = 982.0345499999999
$ws.Range("N136").Value = -12099.1002
